# Trade #25 closed at 2026-02-17 04:08:48 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: refresh aggregate metrics now that trade #25 has closed
#  - Strategy Status sheet: refresh MarketMaking strategy row metrics
#  - All Trades / MarketMaking sheets: mark trade row 26 (Trade #25) as CLOSED
#    with its final P&L / capital / exit data

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Cells.Item(3, 2).Value = 1200.31   # Current Capital
$summary.Cells.Item(4, 2).Value = 0.31      # Total P&L $
$summary.Cells.Item(6, 2).Value = 25        # Total Trades
$summary.Cells.Item(7, 2).Value = 10        # Winning Trades
$summary.Cells.Item(9, 2).Value = 40        # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(4, 3).Value = 100.31     # Capital
$status.Cells.Item(4, 4).Value = 25         # Trades
$status.Cells.Item(4, 5).Value = 0.31       # P&L $
$status.Cells.Item(4, 6).Value = 0.31       # P&L %
$status.Cells.Item(4, 7).Value = 40         # Win Rate %

# ---------------------------------------------------------------------------
# All Trades & MarketMaking sheets - trade row 26 (Trade #25) closes out
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(26, 7).Value = 0.58              # G26 Exit Price
    $ws.Cells.Item(26, 8).Value = "CLOSED"          # H26 Status
    $ws.Cells.Item(26, 9).Value = 1.7544            # I26 P&L %
    $ws.Cells.Item(26, 10).Value = 0.01             # J26 P&L $
    $ws.Cells.Item(26, 11).Value = 100.31           # K26 Capital After
    $ws.Cells.Item(26, 16).Value = "early_exit"     # P26 Exit Reason
    $ws.Cells.Item(26, 17).Value = 0.11             # Q26 Duration (min)
}
